$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("entities")

# Insert a new row *inside* the existing "smartphoneName" merge block (A10:A14)
# so the merged region auto-extends to A10:A15 and every sibling row keeps its
# original formatting. Inserting at row 14 (rather than row 15, the very last
# row of the block) is what makes the engine treat the new row as interior to
# the merge instead of appended after it.
$ws.Rows.Item(14).Insert()

# The insert shifted the old row 14 ("Pixel 2") down to row 15 and left the
# new row 14 blank. Move that content back up to row 14 and put the new
# "G6" entity values into row 15, which is where the new entry belongs.
$ws.Range("B14").Value = $ws.Range("B15").Value2
$ws.Range("C14").Value = $ws.Range("C15").Value2

$ws.Range("B15").Value = "G6"
$ws.Range("C15").Value = "g6"

# Match the selection left behind in the authored workbook.
$ws.Range("A10:A15").Select()

$wb.Save()
